$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Insert two new rows after row 11 (duplicating the "Contact" row) so that
# all IG authors show up as separate Contact rows.
$ws.Rows(12).Insert()
$ws.Rows(12).Insert()

$ws.Range("A12").Value = "Contact"
$ws.Range("B12").Value = "No display for ContactDetail"
$ws.Range("A13").Value = "Contact"
$ws.Range("B13").Value = "No display for ContactDetail"

# Copy styling from the existing contact row so the new rows look the same.
$ws.Range("A11:B11").Copy()
$ws.Range("A12:B13").PasteSpecial(-4122)
$ws.Range("A11:B11").Copy()
$ws.Range("A14:B14").PasteSpecial(-4122)

# Update the generation timestamp cell.
$ws.Range("B8").Value = "2022-01-21T07:49:24+01:00"
